$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row number -> new DAMSLTag (column I) and DialogAct (column J) values
# following a re-run of SGNN to annotate dialog acts.
$updates = @(
    @{ Row = 9; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 11; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 17; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 18; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 32; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 62; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 67; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 73; I = '%'; J = 'Uninterpretable' },
    @{ Row = 74; I = '%'; J = 'Uninterpretable' },
    @{ Row = 77; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 79; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 81; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 95; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 96; I = 'ba'; J = 'Appreciation' },
    @{ Row = 114; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 145; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 153; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 156; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 160; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 175; I = 'qy'; J = 'Yes-No-Question' },
    @{ Row = 178; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 202; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 205; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 212; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 216; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 259; I = '%'; J = 'Uninterpretable' },
    @{ Row = 261; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 263; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 266; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 267; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 268; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 275; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 279; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 288; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 289; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 292; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 300; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 325; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 330; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 373; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 377; I = 'ba'; J = 'Appreciation' },
    @{ Row = 395; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 408; I = '%'; J = 'Uninterpretable' },
    @{ Row = 410; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 433; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 441; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 444; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 472; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 474; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 495; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 500; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 501; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 510; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 515; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 519; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 520; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 523; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 524; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 528; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 532; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 534; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 535; I = 'sv'; J = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
